# Auto-generated script to apply Mandragora_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 15339
$ws.Range("I21").Value = 15339
$ws.Range("K21").Value = 15339
$ws.Range("M21").Value = -14871
$ws.Range("H23").Value = 15339
$ws.Range("I23").Value = 15339
$ws.Range("K23").Value = 15339
$ws.Range("M23").Value = -15105
$ws.Range("H80").Value = 346.45456
$ws.Range("I80").Value = 377.7
$ws.Range("J80").Value = 279.5
$ws.Range("K80").Value = 1133.1
$ws.Range("L80").Value = 838.5
$ws.Range("M80").Value = -135.0999999999999
$ws.Range("N80").Value = -2834.5
$ws.Range("H83").Value = 346.45456
$ws.Range("I83").Value = 377.7
$ws.Range("J83").Value = 279.5
$ws.Range("K83").Value = 3399.3
$ws.Range("L83").Value = 2515.5
$ws.Range("M83").Value = 1592.7
$ws.Range("N83").Value = -12499.5
$ws.Range("H98").Value = 2871.3635
$ws.Range("I98").Value = 3057.9
$ws.Range("J98").Value = 1006
$ws.Range("K98").Value = 3057.9
$ws.Range("L98").Value = 1006
$ws.Range("M98").Value = -1559.9
$ws.Range("N98").Value = -4002
$ws.Range("H122").Value = 2871.3635
$ws.Range("I122").Value = 3057.9
$ws.Range("J122").Value = 1006
$ws.Range("K122").Value = 9173.700000000001
$ws.Range("L122").Value = 3018
$ws.Range("M122").Value = -6723.700000000001
$ws.Range("N122").Value = -7918
$ws.Range("H132").Value = 9636.833000000001
$ws.Range("I132").Value = 9283.611000000001
$ws.Range("K132").Value = 27850.833
$ws.Range("M132").Value = -25320.833
$ws.Range("H138").Value = 22943.89
$ws.Range("I138").Value = 886.0645
$ws.Range("J138").Value = 68530.07000000001
$ws.Range("K138").Value = 2658.1935
$ws.Range("L138").Value = 205590.21
$ws.Range("M138").Value = 2481.8065
$ws.Range("N138").Value = -215870.21

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8349517
$ws.Range("I32").Value = 6862
$ws.Range("J32").Value = 47679176
$ws.Range("K32").Value = 6862
$ws.Range("L32").Value = 47679176
$ws.Range("M32").Value = -6575
$ws.Range("N32").Value = -47679750
$ws.Range("H36").Value = 10000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 10000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 10000
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -10692
$ws.Range("H135").Value = 333353340
$ws.Range("J135").Value = 333353340
$ws.Range("L135").Value = 333353340
$ws.Range("N135").Value = -333363480

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1393.6
$ws.Range("I86").Value = 1338.7693
$ws.Range("K86").Value = 1338.7693
$ws.Range("M86").Value = -215.7692999999999
$ws.Range("H89").Value = 1393.6
$ws.Range("I89").Value = 1338.7693
$ws.Range("K89").Value = 6693.8465
$ws.Range("M89").Value = -1077.8465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 501.75
$ws.Range("I2").Value = 501.75
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 501.75
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -388.75
$ws.Range("N2").ClearContents()
$ws.Range("H31").Value = 1810.2778
$ws.Range("I31").Value = 1417.7
$ws.Range("J31").Value = 2301
$ws.Range("K31").Value = 1417.7
$ws.Range("L31").Value = 2301
$ws.Range("M31").Value = -1122.7
$ws.Range("N31").Value = -2891
$ws.Range("H34").Value = 1810.2778
$ws.Range("I34").Value = 1417.7
$ws.Range("J34").Value = 2301
$ws.Range("K34").Value = 1417.7
$ws.Range("L34").Value = 2301
$ws.Range("M34").Value = -1215.7
$ws.Range("N34").Value = -2705
$ws.Range("H105").Value = 2001400
$ws.Range("I105").Value = 2501250
$ws.Range("K105").Value = 2501250
$ws.Range("M105").Value = -2499503
$ws.Range("H132").Value = 4548241
$ws.Range("I132").Value = 9093391
$ws.Range("J132").Value = 3090.182
$ws.Range("K132").Value = 27280173
$ws.Range("L132").Value = 9270.545999999998
$ws.Range("M132").Value = -27277643
$ws.Range("N132").Value = -14330.546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1400.3334
$ws.Range("J4").Value = 2000
$ws.Range("L4").Value = 6000
$ws.Range("N4").Value = -6224
$ws.Range("H8").Value = 55.214287
$ws.Range("I8").Value = 55.214287
$ws.Range("K8").Value = 165.642861
$ws.Range("M8").Value = -26.64286099999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 1503507.5
$ws.Range("I35").Value = 1503507.5
$ws.Range("K35").Value = 1503507.5
$ws.Range("M35").Value = -1503209.5
$ws.Range("H134").Value = 29065.2
$ws.Range("J134").Value = 29065.2
$ws.Range("L134").Value = 87195.60000000001
$ws.Range("N134").Value = -92265.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3156.182
$ws.Range("I61").Value = 1782.4
$ws.Range("J61").Value = 6100
$ws.Range("K61").Value = 1782.4
$ws.Range("L61").Value = 6100
$ws.Range("M61").Value = -1580.4
$ws.Range("N61").Value = -6504
$ws.Range("H82").Value = 1821
$ws.Range("I82").Value = 1681.5
$ws.Range("J82").Value = 1932.6
$ws.Range("K82").Value = 1681.5
$ws.Range("L82").Value = 1932.6
$ws.Range("M82").Value = -1320.5
$ws.Range("N82").Value = -2654.6
$ws.Range("H85").Value = 1821
$ws.Range("I85").Value = 1681.5
$ws.Range("J85").Value = 1932.6
$ws.Range("K85").Value = 1681.5
$ws.Range("L85").Value = 1932.6
$ws.Range("M85").Value = -433.5
$ws.Range("N85").Value = -4428.6
$ws.Range("H113").Value = 3156.182
$ws.Range("I113").Value = 1782.4
$ws.Range("J113").Value = 6100
$ws.Range("K113").Value = 1782.4
$ws.Range("L113").Value = 6100
$ws.Range("M113").Value = 387.5999999999999
$ws.Range("N113").Value = -10440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 3533.3333
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 4800
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 4800
$ws.Range("M8").Value = -860
$ws.Range("N8").Value = -5080
$ws.Range("H30").Value = 50006.668
$ws.Range("J30").Value = 50006.668
$ws.Range("L30").Value = 50006.668
$ws.Range("N30").Value = -50220.668
$ws.Range("H81").Value = 1011.0833
$ws.Range("I81").Value = 731
$ws.Range("J81").Value = 1211.1428
$ws.Range("K81").Value = 1462
$ws.Range("L81").Value = 2422.2856
$ws.Range("M81").Value = -401
$ws.Range("N81").Value = -4544.2856
$ws.Range("H84").Value = 1011.0833
$ws.Range("I84").Value = 731
$ws.Range("J84").Value = 1211.1428
$ws.Range("K84").Value = 7310
$ws.Range("L84").Value = 12111.428
$ws.Range("M84").Value = -2006
$ws.Range("N84").Value = -22719.428
$ws.Range("H108").Value = 28980
$ws.Range("J108").Value = 28980
$ws.Range("L108").Value = 28980
$ws.Range("N108").Value = -36660
